$wb = $excel.ActiveWorkbook
$excel.Calculation = -4135  # xlCalculationManual

# ---- Sheet: FS ----
$ws = $wb.Worksheets.Item("FS")
$ws.Range("C18").Value = 1.0
$ws.Range("D18").Value = 0.0
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 31.57894736842105
$ws.Range("J18").Value = 0.9333333333333333
$ws.Range("K18").Value = 0.06666666666666667
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 11.695906432748536
$ws.Range("C19").Value = 0.125
$ws.Range("D19").Value = 0.625
$ws.Range("E19").Value = 0.25
$ws.Range("J19").Value = 0.05555555555555555
$ws.Range("K19").Value = 0.8333333333333334
$ws.Range("L19").Value = 0.1111111111111111
$ws.Range("C20").Value = 0.16666666666666666
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.5
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.09259259259259259
$ws.Range("L20").Value = 0.9074074074074074
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.8
$ws.Range("K24").Value = 0.2
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 29.239766081871345
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.625
$ws.Range("E25").Value = 0.375
$ws.Range("J25").Value = 0.08333333333333333
$ws.Range("K25").Value = 0.6666666666666666
$ws.Range("L25").Value = 0.25
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.8333333333333334
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 0.3148148148148148
$ws.Range("L26").Value = 0.6851851851851852

# ---- Sheet: IF ----
$ws = $wb.Worksheets.Item("IF")
$ws.Range("C18").Value = 0.2
$ws.Range("D18").Value = 0.6
$ws.Range("E18").Value = 0.2
$ws.Range("G18").Value = 47.368421052631575
$ws.Range("J18").Value = 0.6
$ws.Range("K18").Value = 0.35555555555555557
$ws.Range("L18").Value = 0.044444444444444446
$ws.Range("N18").Value = 33.91812865497076
$ws.Range("C19").Value = 0.375
$ws.Range("D19").Value = 0.625
$ws.Range("E19").Value = 0.0
$ws.Range("J19").Value = 0.19444444444444445
$ws.Range("K19").Value = 0.6527777777777778
$ws.Range("L19").Value = 0.1527777777777778
$ws.Range("C20").Value = 0.16666666666666666
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("E20").Value = 0.6666666666666666
$ws.Range("J20").Value = 0.018518518518518517
$ws.Range("K20").Value = 0.25925925925925924
$ws.Range("L20").Value = 0.7222222222222222
$ws.Range("C24").Value = 0.6
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.4
$ws.Range("G24").Value = 42.10526315789473
$ws.Range("J24").Value = 0.7555555555555555
$ws.Range("K24").Value = 0.15555555555555556
$ws.Range("L24").Value = 0.08888888888888889
$ws.Range("N24").Value = 27.485380116959064
$ws.Range("C25").Value = 0.125
$ws.Range("D25").Value = 0.625
$ws.Range("E25").Value = 0.25
$ws.Range("J25").Value = 0.06944444444444445
$ws.Range("K25").Value = 0.7361111111111112
$ws.Range("L25").Value = 0.19444444444444445
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.5
$ws.Range("E26").Value = 0.5
$ws.Range("J26").Value = 0.09259259259259259
$ws.Range("K26").Value = 0.2222222222222222
$ws.Range("L26").Value = 0.6851851851851852

# ---- Sheet: IA ----
$ws = $wb.Worksheets.Item("IA")
$ws.Range("C18").Value = 1.0
$ws.Range("D18").Value = 0.0
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 15.789473684210526
$ws.Range("J18").Value = 0.9333333333333333
$ws.Range("K18").Value = 0.06666666666666667
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 11.11111111111111
$ws.Range("C19").Value = 0.25
$ws.Range("D19").Value = 0.75
$ws.Range("E19").Value = 0.0
$ws.Range("J19").Value = 0.05555555555555555
$ws.Range("K19").Value = 0.8611111111111112
$ws.Range("L19").Value = 0.08333333333333333
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("E20").Value = 0.8333333333333334
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.1111111111111111
$ws.Range("L20").Value = 0.8888888888888888
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 31.57894736842105
$ws.Range("J24").Value = 0.8222222222222222
$ws.Range("K24").Value = 0.08888888888888889
$ws.Range("L24").Value = 0.08888888888888889
$ws.Range("N24").Value = 26.900584795321635
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 0.75
$ws.Range("E25").Value = 0.25
$ws.Range("J25").Value = 0.041666666666666664
$ws.Range("K25").Value = 0.75
$ws.Range("L25").Value = 0.20833333333333334
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.6666666666666666
$ws.Range("E26").Value = 0.3333333333333333
$ws.Range("J26").Value = 0.09259259259259259
$ws.Range("K26").Value = 0.2777777777777778
$ws.Range("L26").Value = 0.6296296296296297

# ---- Sheet: FS-IF ----
$ws = $wb.Worksheets.Item("FS-IF")
$ws.Range("C18").Value = 0.8
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 15.789473684210526
$ws.Range("J18").Value = 0.6444444444444445
$ws.Range("K18").Value = 0.3333333333333333
$ws.Range("L18").Value = 0.022222222222222223
$ws.Range("N18").Value = 23.391812865497073
$ws.Range("C19").Value = 0.0
$ws.Range("D19").Value = 0.875
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.1111111111111111
$ws.Range("K19").Value = 0.7777777777777778
$ws.Range("L19").Value = 0.1111111111111111
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("E20").Value = 0.8333333333333334
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.14814814814814814
$ws.Range("L20").Value = 0.8518518518518519
$ws.Range("C24").Value = 0.6
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.2
$ws.Range("G24").Value = 21.052631578947366
$ws.Range("J24").Value = 0.8888888888888888
$ws.Range("K24").Value = 0.1111111111111111
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 23.391812865497073
$ws.Range("C25").Value = 0.125
$ws.Range("D25").Value = 0.75
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.06944444444444445
$ws.Range("K25").Value = 0.7361111111111112
$ws.Range("L25").Value = 0.19444444444444445
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.0
$ws.Range("E26").Value = 1.0
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.25925925925925924
$ws.Range("L26").Value = 0.7037037037037037

# ---- Sheet: FS-IA ----
$ws = $wb.Worksheets.Item("FS-IA")
$ws.Range("C18").Value = 0.8
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 31.57894736842105
$ws.Range("J18").Value = 0.8444444444444444
$ws.Range("K18").Value = 0.15555555555555556
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 14.035087719298245
$ws.Range("C19").Value = 0.25
$ws.Range("D19").Value = 0.625
$ws.Range("E19").Value = 0.125
$ws.Range("J19").Value = 0.027777777777777776
$ws.Range("K19").Value = 0.8611111111111112
$ws.Range("L19").Value = 0.1111111111111111
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.6666666666666666
$ws.Range("J20").Value = 0.018518518518518517
$ws.Range("K20").Value = 0.1111111111111111
$ws.Range("L20").Value = 0.8703703703703703
$ws.Range("C24").Value = 1.0
$ws.Range("D24").Value = 0.0
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 26.31578947368421
$ws.Range("J24").Value = 0.8666666666666667
$ws.Range("K24").Value = 0.1111111111111111
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 23.391812865497073
$ws.Range("C25").Value = 0.375
$ws.Range("D25").Value = 0.5
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.027777777777777776
$ws.Range("K25").Value = 0.75
$ws.Range("L25").Value = 0.2222222222222222
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.16666666666666666
$ws.Range("E26").Value = 0.8333333333333334
$ws.Range("J26").Value = 0.037037037037037035
$ws.Range("K26").Value = 0.25925925925925924
$ws.Range("L26").Value = 0.7037037037037037

# ---- Sheet: IF-IA ----
$ws = $wb.Worksheets.Item("IF-IA")
$ws.Range("C18").Value = 0.8
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 21.052631578947366
$ws.Range("J18").Value = 0.6222222222222222
$ws.Range("K18").Value = 0.3333333333333333
$ws.Range("L18").Value = 0.044444444444444446
$ws.Range("N18").Value = 24.561403508771928
$ws.Range("C19").Value = 0.125
$ws.Range("D19").Value = 0.875
$ws.Range("E19").Value = 0.0
$ws.Range("J19").Value = 0.1111111111111111
$ws.Range("K19").Value = 0.7777777777777778
$ws.Range("L19").Value = 0.1111111111111111
$ws.Range("C20").Value = 0.16666666666666666
$ws.Range("D20").Value = 0.16666666666666666
$ws.Range("E20").Value = 0.6666666666666666
$ws.Range("J20").Value = 0.018518518518518517
$ws.Range("K20").Value = 0.14814814814814814
$ws.Range("L20").Value = 0.8333333333333334
$ws.Range("C24").Value = 0.8
$ws.Range("D24").Value = 0.2
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 15.789473684210526
$ws.Range("J24").Value = 0.8444444444444444
$ws.Range("K24").Value = 0.13333333333333333
$ws.Range("L24").Value = 0.022222222222222223
$ws.Range("N24").Value = 21.637426900584796
$ws.Range("C25").Value = 0.0
$ws.Range("D25").Value = 1.0
$ws.Range("E25").Value = 0.0
$ws.Range("J25").Value = 0.05555555555555555
$ws.Range("K25").Value = 0.7638888888888888
$ws.Range("L25").Value = 0.18055555555555555
$ws.Range("C26").Value = 0.0
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.6666666666666666
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 0.24074074074074073
$ws.Range("L26").Value = 0.7592592592592593

# ---- Sheet: FS-IF-IA ----
$ws = $wb.Worksheets.Item("FS-IF-IA")
$ws.Range("C12").Value = 0.6
$ws.Range("D12").Value = 0.4
$ws.Range("E12").Value = 0.0
$ws.Range("G12").Value = 26.31578947368421
$ws.Range("J12").Value = 0.8666666666666667
$ws.Range("K12").Value = 0.044444444444444446
$ws.Range("L12").Value = 0.08888888888888889
$ws.Range("N12").Value = 16.374269005847953
$ws.Range("C13").Value = 0.125
$ws.Range("D13").Value = 0.875
$ws.Range("E13").Value = 0.0
$ws.Range("J13").Value = 0.041666666666666664
$ws.Range("K13").Value = 0.8611111111111112
$ws.Range("L13").Value = 0.09722222222222222
$ws.Range("C14").Value = 0.0
$ws.Range("D14").Value = 0.3333333333333333
$ws.Range("E14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.07407407407407407
$ws.Range("K14").Value = 0.14814814814814814
$ws.Range("L14").Value = 0.7777777777777778
$ws.Range("C18").Value = 0.8
$ws.Range("D18").Value = 0.2
$ws.Range("E18").Value = 0.0
$ws.Range("G18").Value = 21.052631578947366
$ws.Range("J18").Value = 0.7777777777777778
$ws.Range("K18").Value = 0.2222222222222222
$ws.Range("L18").Value = 0.0
$ws.Range("N18").Value = 19.883040935672515
$ws.Range("C19").Value = 0.125
$ws.Range("D19").Value = 0.875
$ws.Range("E19").Value = 0.0
$ws.Range("J19").Value = 0.1388888888888889
$ws.Range("K19").Value = 0.7361111111111112
$ws.Range("L19").Value = 0.125
$ws.Range("C20").Value = 0.0
$ws.Range("D20").Value = 0.3333333333333333
$ws.Range("E20").Value = 0.6666666666666666
$ws.Range("J20").Value = 0.0
$ws.Range("K20").Value = 0.09259259259259259
$ws.Range("L20").Value = 0.9074074074074074
$ws.Range("C24").Value = 0.6
$ws.Range("D24").Value = 0.4
$ws.Range("E24").Value = 0.0
$ws.Range("G24").Value = 36.84210526315789
$ws.Range("J24").Value = 0.8666666666666667
$ws.Range("K24").Value = 0.13333333333333333
$ws.Range("L24").Value = 0.0
$ws.Range("N24").Value = 23.391812865497073
$ws.Range("C25").Value = 0.125
$ws.Range("D25").Value = 0.75
$ws.Range("E25").Value = 0.125
$ws.Range("J25").Value = 0.05555555555555555
$ws.Range("K25").Value = 0.7361111111111112
$ws.Range("L25").Value = 0.20833333333333334
$ws.Range("C26").Value = 0.16666666666666666
$ws.Range("D26").Value = 0.3333333333333333
$ws.Range("E26").Value = 0.5
$ws.Range("J26").Value = 0.0
$ws.Range("K26").Value = 0.2777777777777778
$ws.Range("L26").Value = 0.7222222222222222
